$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.328.03'
$ws.Range("E2").Value = '  -0.06%  '
$ws.Range("D3").Value = '3.904.47'
$ws.Range("E3").Value = '  -0.12%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '487.22'
$ws.Range("E5").Value = '  +1.41%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.02'
$ws.Range("E6").Value = '  +0.82%  '
$ws.Range("E7").Value = '  +0.34%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.744'
$ws.Range("E9").Value = '  +2.69%  '
$ws.Range("E10").Value = '  +8.61%  '
$ws.Range("E11").Value = '  +1.74%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '43.21'
$ws.Range("E12").Value = '  +1.29%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '10.49'
$ws.Range("E13").Value = '  -1.46%  '
$ws.Range("D14").Value = '4.519.33'
$ws.Range("E14").Value = '  -0.23%  '
$ws.Range("D15").Value = '3.902.29'
$ws.Range("E15").Value = '  -0.91%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.23'
$ws.Range("E16").Value = '  -2.67%  '
$ws.Range("E17").Value = '  -0.49%  '
$ws.Range("E18").Value = '  +1.54%  '
$ws.Range("E19").Value = '  +1.78%  '
$ws.Range("D20").Value = '68.394.11'
$ws.Range("E20").Value = '  +0.04%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '433.18'
$ws.Range("E21").Value = '  -0.51%  '
$ws.Range("E22").Value = '  +5.63%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '14.71'
$ws.Range("E23").Value = '  +0.30%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '89.85'
$ws.Range("E24").Value = '  +2.33%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.35'
$ws.Range("E25").Value = '  +18.45%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.75'
$ws.Range("E26").Value = '  +4.49%  '
$ws.Range("E27").Value = '  -6.29%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '37.53'
$ws.Range("E28").Value = '  -1.57%  '
$ws.Range("E29").Value = '  -2.19%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '712.93'
$ws.Range("E30").Value = '  +0.90%  '
$ws.Range("E31").Value = '  +0.69%  '
$ws.Range("E32").Value = '  +0.41%  '
$ws.Range("E33").Value = '  +2.13%  '
$ws.Range("E34").Value = '  -3.51%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.09'
$ws.Range("E35").Value = '  +6.14%  '
$ws.Range("E36").Value = '  +3.61%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '40.84'
$ws.Range("E37").Value = '  -1.82%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.405'
$ws.Range("E38").Value = '  +18.85%  '
$ws.Range("E39").Value = '  -1.91%  '
$ws.Range("E40").Value = '  +0.05%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0495'
$ws.Range("E41").Value = '  +4.65%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.96'
$ws.Range("E42").Value = '  +7.56%  '
$ws.Range("E43").Value = '  +3.06%  '
$ws.Range("E44").Value = '  -1.15%  '
$ws.Range("D45").Value = '0.0₆0380'
$ws.Range("E45").Value = '  +29.30%  '
$ws.Range("E46").Value = '  +1.14%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.38'
$ws.Range("E47").Value = '  +7.35%  '
$ws.Range("E48").Value = '  +0.22%  '
$ws.Range("E49").Value = '  -1.59%  '
$ws.Range("E50").Value = '  -1.99%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '142.81'
$ws.Range("E51").Value = '  -2.20%  '
